$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = "39.195.86"
$ws.Cells.Item(2, 5).Value = "  -2.17%  "
$ws.Cells.Item(3, 4).Value = "2.196.77"
$ws.Cells.Item(3, 5).Value = "  -5.83%  "
$ws.Cells.Item(4, 5).Value = "  +0.08%  "
$ws.Cells.Item(5, 4).Value = "'294.18"
$ws.Cells.Item(5, 5).Value = "  -4.29%  "
$ws.Cells.Item(6, 4).Value = "'80.80"
$ws.Cells.Item(6, 5).Value = "  -4.91%  "
$ws.Cells.Item(7, 4).Value = "'0.509"
$ws.Cells.Item(7, 5).Value = "  -3.92%  "
$ws.Cells.Item(8, 5).Value = "  +0.03%  "
$ws.Cells.Item(9, 4).Value = "'0.465"
$ws.Cells.Item(9, 5).Value = "  -3.70%  "
$ws.Cells.Item(10, 4).Value = "'0.0767"
$ws.Cells.Item(10, 5).Value = "  -6.02%  "
$ws.Cells.Item(11, 4).Value = "'28.96"
$ws.Cells.Item(11, 5).Value = "  -3.75%  "
$ws.Cells.Item(12, 4).Value = "'46.74"
$ws.Cells.Item(12, 5).Value = "  -11.23%  "
$ws.Cells.Item(13, 5).Value = "  -2.72%  "
$ws.Cells.Item(14, 4).Value = "2.532.86"
$ws.Cells.Item(14, 5).Value = "  -5.93%  "
$ws.Cells.Item(15, 4).Value = "'6.20"
$ws.Cells.Item(15, 5).Value = "  -3.21%  "
$ws.Cells.Item(16, 4).Value = "'13.86"
$ws.Cells.Item(16, 5).Value = "  -5.62%  "
$ws.Cells.Item(17, 4).Value = "2.194.81"
$ws.Cells.Item(17, 5).Value = "  -5.60%  "
$ws.Cells.Item(18, 5).Value = "  -5.90%  "
$ws.Cells.Item(19, 4).Value = "39.096.21"
$ws.Cells.Item(19, 5).Value = "  -2.35%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0868"
$ws.Cells.Item(20, 5).Value = "  -3.84%  "
$ws.Cells.Item(21, 4).Value = "'5.68"
$ws.Cells.Item(21, 5).Value = "  -6.63%  "
$ws.Cells.Item(22, 4).Value = "'64.45"
$ws.Cells.Item(22, 5).Value = "  -4.62%  "
$ws.Cells.Item(23, 5).Value = "  -4.18%  "
$ws.Cells.Item(24, 4).Value = "'225.25"
$ws.Cells.Item(24, 5).Value = "  -4.31%  "
$ws.Cells.Item(25, 5).Value = "  -0.16%  "
$ws.Cells.Item(26, 4).Value = "'2.39"
$ws.Cells.Item(26, 5).Value = "  -6.79%  "
$ws.Cells.Item(27, 5).Value = "  -0.70%  "
$ws.Cells.Item(28, 4).Value = "'22.41"
$ws.Cells.Item(28, 5).Value = "  -4.07%  "
$ws.Cells.Item(29, 5).Value = "  -2.03%  "
$ws.Cells.Item(30, 4).Value = "'9.01"
$ws.Cells.Item(30, 5).Value = "  -2.80%  "
$ws.Cells.Item(31, 4).Value = "'149.22"
$ws.Cells.Item(31, 5).Value = "  -1.79%  "
$ws.Cells.Item(32, 4).Value = "'31.45"
$ws.Cells.Item(32, 5).Value = "  -11.22%  "
$ws.Cells.Item(33, 5).Value = "  -0.16%  "
$ws.Cells.Item(34, 4).Value = "'4.77"
$ws.Cells.Item(34, 5).Value = "  -6.60%  "
$ws.Cells.Item(35, 5).Value = "  -4.32%  "
$ws.Cells.Item(36, 4).Value = "'0.0692"
$ws.Cells.Item(36, 5).Value = "  -4.33%  "
$ws.Cells.Item(37, 5).Value = "  -3.72%  "
$ws.Cells.Item(38, 4).Value = "'15.21"
$ws.Cells.Item(38, 5).Value = "  -3.47%  "
$ws.Cells.Item(39, 4).Value = "'0.0953"
$ws.Cells.Item(39, 5).Value = "  -4.42%  "
$ws.Cells.Item(40, 5).Value = "  -5.42%  "
$ws.Cells.Item(41, 4).Value = "'1.64"
$ws.Cells.Item(41, 5).Value = "  -3.72%  "
$ws.Cells.Item(42, 4).Value = "'3.59"
$ws.Cells.Item(42, 5).Value = "  -5.86%  "
$ws.Cells.Item(43, 4).Value = "1.897.14"
$ws.Cells.Item(43, 5).Value = "  -2.18%  "
$ws.Cells.Item(44, 4).Value = "'2.05"
$ws.Cells.Item(44, 5).Value = "  -9.66%  "
$ws.Cells.Item(45, 4).Value = "'0.0259"
$ws.Cells.Item(45, 5).Value = "  -3.00%  "
$ws.Cells.Item(46, 4).Value = "'8.96"
$ws.Cells.Item(46, 5).Value = "  -3.23%  "
$ws.Cells.Item(47, 4).Value = "'15.84"
$ws.Cells.Item(47, 5).Value = "  -10.20%  "
$ws.Cells.Item(48, 4).Value = "'2.59"
$ws.Cells.Item(48, 5).Value = "  -3.26%  "
$ws.Cells.Item(49, 4).Value = "2.404.24"
$ws.Cells.Item(49, 5).Value = "  -6.04%  "
$ws.Cells.Item(50, 4).Value = "'70.92"
$ws.Cells.Item(50, 5).Value = "  -0.80%  "
$ws.Cells.Item(51, 4).Value = "'86.88"
$ws.Cells.Item(51, 5).Value = "  -6.47%  "
